$d = $word.ActiveDocument

$replacements = @(
    @{old = "100÷3="; new = "799÷8="},
    @{old = "588÷2="; new = "521÷4="},
    @{old = "171÷3="; new = "767÷4="},
    @{old = "118÷9="; new = "904÷9="},
    @{old = "107÷7="; new = "726÷9="},
    @{old = "781÷5="; new = "661÷4="},
    @{old = "219÷5="; new = "995÷3="},
    @{old = "232÷9="; new = "123÷4="},
    @{old = "567÷3="; new = "853÷4="},
    @{old = "308÷4="; new = "448÷3="},
    @{old = "779÷9="; new = "303÷3="},
    @{old = "653÷7="; new = "208÷4="},
    @{old = "446÷2="; new = "426÷9="},
    @{old = "731÷6="; new = "132÷4="},
    @{old = "614÷9="; new = "401÷8="},
    @{old = "586÷8="; new = "490÷9="},
    @{old = "174÷8="; new = "273÷9="},
    @{old = "661÷5="; new = "859÷7="},
    @{old = "759÷4="; new = "467÷8="},
    @{old = "507÷8="; new = "854÷5="},
    @{old = "308÷2="; new = "731÷2="},
    @{old = "572÷7="; new = "941÷3="},
    @{old = "759÷6="; new = "961÷4="},
    @{old = "605÷9="; new = "380÷4="},
    @{old = "120÷7="; new = "698÷7="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
